$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row labels: "<field>_old" -> "<field>_FV2310", "<field>_new" -> "<field>_FV2404"
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# Freeze header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a real Excel table (ListObject)
$range = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
